# Update gh-pages to output generated at 456a3b4
#
# Both the "展览" (exhibition) sheet (index 1) and the "全部类型" (all
# types) sheet (index 4) share the same first 21 data rows. A new event
# ("合肥·人气COSER次元会（免费展）", 2024-06-16) is inserted as row 7,
# pushing every following row down by one. A few "想去人数" (interest
# count) values were refreshed at the same time, and two totals above the
# insertion point (F2, F4) were bumped as well.

$wb = $excel.ActiveWorkbook

# Sheets that contain the affected table (by 1-based index):
#   1 = 展览 (Exhibition)
#   4 = 全部类型 (All types)
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # --- Refresh the two "want to go" counters above the insertion point ---
    $ws.Cells.Item(2, 6).Value = 321
    $ws.Cells.Item(4, 6).Value = 10388

    # --- Insert the new row at position 7, pushing rows 7.. down by one ---
    $ws.Rows.Item(7).Insert()

    # New row 7's index cell (column A) needs to match the bold/centered/
    # bordered style used by every other index cell in column A, since a
    # plain row insert leaves it with a border-less copy of that style.
    $a7 = $ws.Cells.Item(7, 1)
    $a7.Font.Bold = $true
    $a7.HorizontalAlignment = -4108
    $a7.VerticalAlignment = -4160
    $a7.Borders.LineStyle = 1
    $a7.Value = 6

    # Column B holds plain-text "YYYY-MM-DD" strings. Excel's COM layer
    # auto-converts that pattern into a real date serial, so force the
    # cell to Text first, then reset the visual style back to Normal so
    # no stray NumberFormat/style id lingers on the cell once saved.
    $b7 = $ws.Cells.Item(7, 2)
    $b7.NumberFormat = "@"
    $b7.Value = "2024-06-16"
    $b7.Style = "Normal"

    $ws.Cells.Item(7, 3).Value = "合肥·人气COSER次元会（免费展）"
    $ws.Cells.Item(7, 4).Value = "长江中路98号 合肥银泰中心"
    $ws.Cells.Item(7, 5).Value = "2024.06.16 12:30-06.16 20:00"
    $ws.Cells.Item(7, 6).Value = 10
    $ws.Cells.Item(7, 7).Value = 29.9
    $ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86773"
    $ws.Cells.Item(7, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/wbdsV9SD1717387577908.png"

    # --- Refresh "想去人数" (want-to-go count) on the rows that shifted down ---
    $ws.Cells.Item(8, 6).Value = 1283
    $ws.Cells.Item(9, 6).Value = 7189
    $ws.Cells.Item(10, 6).Value = 18
    $ws.Cells.Item(11, 6).Value = 442
    $ws.Cells.Item(12, 6).Value = 205
    $ws.Cells.Item(13, 6).Value = 127
    $ws.Cells.Item(14, 6).Value = 3205
    $ws.Cells.Item(15, 6).Value = 37
    $ws.Cells.Item(16, 6).Value = 314
    $ws.Cells.Item(17, 6).Value = 665
    $ws.Cells.Item(18, 6).Value = 123
    $ws.Cells.Item(19, 6).Value = 1037
    $ws.Cells.Item(20, 6).Value = 278
    $ws.Cells.Item(21, 6).Value = 73
    $ws.Cells.Item(22, 6).Value = 1640

    # Column A is a plain sequential row index (row number - 1, header is
    # 0) — it is NOT event data, so it was never meant to shift along with
    # the rest of the row when the new row was inserted. A plain row
    # Insert() drags the old index value down with the rest of the row,
    # so re-number A8..A22 back to the correct sequential values.
    for ($r = 8; $r -le 22; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}
